$d = $word.ActiveDocument

# The document embeds two logos (Pearson "image1.png" and BTEC "image2.jpg")
# twice each (in the "default" and "first page" header/footer parts). Each
# picture's drawing carries its filename in two places in the OOXML:
#   <wp:docPr .../>      (the drawing's docPr "name" attribute)
#   <pic:cNvPr .../>     (the picture shape's "name" attribute, nested inside)
# The edit swaps the logical file names used in those "name" attributes
# (Pearson logo: image1.png -> image2.png; BTEC logo: image2.jpg -> image1.jpg)
# while leaving the "descr", "id" attributes and the actual media parts/
# relationships untouched.
#
# The Word object model doesn't expose a settable "name" property on
# InlineShape that reaches the nested <pic:cNvPr> element, so we perform the
# swap directly on the document's flat OPC XML (WordOpenXML), which is a
# faithful, round-trippable textual representation of the package.

$xml = $d.WordOpenXML

function Rename-DrawingName {
    param(
        [string]$Text,
        [string]$DocPrId,
        [string]$OldName,
        [string]$NewName
    )

    # Unique anchor for this particular drawing's docPr element.
    $docPrNeedle = 'id="' + $DocPrId + '" name="' + $OldName + '"/>'
    $docPrReplacement = 'id="' + $DocPrId + '" name="' + $NewName + '"/>'

    $anchorIndex = $Text.IndexOf($docPrNeedle)
    if ($anchorIndex -lt 0) {
        throw "Could not find wp:docPr anchor for id=$DocPrId name=$OldName"
    }

    # Replace the docPr's name (first hit only, right after the anchor point).
    $before = $Text.Substring(0, $anchorIndex)
    $after = $Text.Substring($anchorIndex)
    $after = $after.Substring(0, $docPrNeedle.Length).Replace($docPrNeedle, $docPrReplacement) + $after.Substring($docPrNeedle.Length)
    $Text = $before + $after

    # Now find the nested <pic:cNvPr ... id="0" name="OldName"/> that belongs
    # to this same drawing -- it is the first such attribute pair following
    # the docPr element we just updated.
    $searchFrom = $anchorIndex
    $cNvNeedle = 'id="0" name="' + $OldName + '"/>'
    $cNvReplacement = 'id="0" name="' + $NewName + '"/>'
    $cNvIndex = $Text.IndexOf($cNvNeedle, $searchFrom)
    if ($cNvIndex -lt 0) {
        throw "Could not find pic:cNvPr for docPr id=$DocPrId name=$OldName"
    }
    $before = $Text.Substring(0, $cNvIndex)
    $after = $Text.Substring($cNvIndex)
    $after = $after.Substring(0, $cNvNeedle.Length).Replace($cNvNeedle, $cNvReplacement) + $after.Substring($cNvNeedle.Length)
    $Text = $before + $after

    return $Text
}

# Pearson logo drawings: id="2" (footer, first page) and id="4" (footer, default)
$xml = Rename-DrawingName $xml "2" "image1.png" "image2.png"
$xml = Rename-DrawingName $xml "4" "image1.png" "image2.png"

# BTEC logo drawings: id="1" (header, first page) and id="3" (header, default)
$xml = Rename-DrawingName $xml "1" "image2.jpg" "image1.jpg"
$xml = Rename-DrawingName $xml "3" "image2.jpg" "image1.jpg"

$d.WordOpenXML = $xml

Write-Host "Renamed drawing name attributes for 4 pictures"
